# Update the "department" column (C) on the courses sheet.
# Previously every row used the single generic label
# "FACULTY OF BUSINESS & TECHNOLOGY"; split it out into the
# specific faculty/stream that each row actually belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Rows 2-5: Business qualifications (TAS, BSB.. courses)
$ws.Range("C2:C5").Value = "Business"

# Rows 6-8: Information Technology qualifications (ICT.. courses)
$ws.Range("C6:C8").Value = "Information Technology"

# Row 9: Building and Construction (RII.. course)
$ws.Range("C9:C9").Value = "Building and Construction"

# Rows 10-15: combined course "Packages"
$ws.Range("C10:C15").Value = "Packages"
